$wb = $excel.ActiveWorkbook
$summary = $wb.Worksheets.Item("Summary")
$new = $wb.Worksheets.Add($summary, [System.Type]::Missing)
$new.Name = "Runtime-Monitor_Patch"
Write-Host $wb.FullName
Write-Host $wb.Path
